# Apply the authored change set to the active presentation:
#  1. Re-apply the table style on the three data tables (slides 14-16)
#     from the custom "Table_0" style to the built-in table style
#     {4B4AC84C-A9E1-4982-87B3-C163C3FEF339}.
#  2. Swap the deck's applied theme colour palette from "Integral"
#     (Red Violet) to the classic "Office" palette - i.e. the colours
#     that the slide master / presentation theme (theme2.xml) resolves
#     to are updated to the Office Theme colour values.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyleId = "{4B4AC84C-A9E1-4982-87B3-C163C3FEF339}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colour palette ------------------------------------------
# RGB() packs colours as 0x00BBGGRR, matching the VBA RGB() convention
# used by ThemeColorScheme.Item(n).RGB.
$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
